$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-slot labels in column C (B column already holds the
# alternating "2x rows" pattern and is untouched).
$ws.Range("C2").Value = "2:55-3:0"
$ws.Range("C3").Value = "3:0-3:5"
$ws.Range("C6").Value = "18:55-19:0"
$ws.Range("C7").Value = "19:0-19:5"

# Move the active selection from C11 to B11.
$ws.Range("B11").Select()
